$d = $word.ActiveDocument

# Locate the list item that starts with "Visualization" (the
# "Visualization - create charts ... relevant queries." bullet) and
# colour its whole paragraph (including the paragraph mark, so the
# bullet/number picks up the colour too) red, to flag it as now
# working per the commit note ("jdbc - chart wroking").
$rng = $d.Content
$found = $rng.Find.Execute("Visualization*charts*relevant queries*", $false, $false, $true,
                            $false, $false, $true, 1, $false, "", 0)

if ($found -and $rng.Find.Found) {
    $target = $rng.Paragraphs(1)
    $target.Range.Font.Color = 255
} else {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "Visualization*charts*relevant queries*") {
            $p.Range.Font.Color = 255
            break
        }
    }
}
